# ProjectConfiguration.xlsx - update the "modelFolder" reference and the
# active selection, matching the upstream commit that wires the config
# sheet up to point at the "../Models" folder instead of "../..".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 is the "modelFolder" property; change its Value column (B2)
# from "../.." to "../Models".
$ws.Range("B2").Value = "../Models"

# The author's last selection before saving was cell B2 (previously B14).
$ws.Range("B2").Select()

$wb.Save()
